# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the
# "© 2020 . Contact: ..." paragraph (and the blank paragraph that separated
# them from the preceding "LOB1012: ..." line), leaving a single blank
# paragraph before the trailing page-break paragraph.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph.
$start = $d.Content
$start.Find.ClearFormatting()
$found = $start.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the "(c) 2020 . Contact: ..." paragraph.
$end = $d.Content
$end.Find.ClearFormatting()
$end.Find.Execute([char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Build a range spanning the start of the blank paragraph immediately before
# "Ver no Jupiter ..." through the end of the "(c) 2020 ..." paragraph
# (including its paragraph mark), then delete it.
$blank = $start.Paragraphs(1).Previous(1)
$r = $d.Range($blank.Range.Start, $end.Paragraphs(1).Range.End)
$r.Delete()
